# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    existing "2022-Q2" sheet) and populate it with the quarter's fund
#    holdings table (same B:H column layout as the other quarter sheets).
# 2) Insert a new summary row into "总计" for "2022-Q3" (count=28,
#    value=0.93) right above the existing "2022-Q2" row, shifting every
#    other row down by one and renumbering the leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Header row (columns B:H - column A is left blank, matching every other
# quarter sheet in this workbook).
$q3Sheet.Cells.Item(1, 2).Value = "基金代码"
$q3Sheet.Cells.Item(1, 3).Value = "基金名称"
$q3Sheet.Cells.Item(1, 4).Value = "基金规模"
$q3Sheet.Cells.Item(1, 5).Value = "股票总仓位"
$q3Sheet.Cells.Item(1, 6).Value = "仓位占比"
$q3Sheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3Sheet.Cells.Item(1, 8).Value = "仓位排名"
$q3Sheet.Range("B1:H1").Font.Bold = $true

# Data rows: code, name, scale, position, ratio, marketValue, rank.
# scale/position/ratio/marketValue are kept as TEXT (they are rendered with
# fixed decimals such as "7.50"/"0.2122" in the source data), so the target
# columns are pre-formatted as Text before the values are written - this
# stops them being auto-coerced into numbers (and losing trailing zeros).
$q3Data = @(
  @('003713', '英大睿盛灵活配置混合A', '2.83', '93.65', '7.50', '0.2122', 2),
  @('003714', '英大睿盛灵活配置混合C', '2.19', '93.65', '7.50', '0.1642', 2),
  @('900009', '中信证券成长动力混合A', '4.01', '88.73', '3.29', '0.1319', 6),
  @('004809', '新疆前海联合润丰灵活配置混合A', '1.30', '87.05', '7.28', '0.0946', 4),
  @('005310', '广发电子信息传媒股票A', '1.55', '89.36', '2.53', '0.0392', 10),
  @('001607', '英大策略优选混合A', '0.57', '91.98', '6.33', '0.0361', 6),
  @('000458', '英大领先回报混合', '1.32', '93.57', '2.16', '0.0285', 4),
  @('004890', '中邮健康文娱灵活配置混合', '0.41', '86.15', '6.11', '0.0251', 2),
  @('012522', '英大稳固增强核心一年持有混合C', '1.24', '27.71', '1.96', '0.0243', 4),
  @('008848', '中融智选对冲策略3个月定期开放灵活配置混合', '2.28', '67.51', '0.95', '0.0217', 8),
  @('900059', '中信证券成长动力混合C', '0.52', '88.73', '3.29', '0.0171', 6),
  @('003447', '英大睿鑫灵活配置混合C', '0.21', '92.71', '7.38', '0.0155', 9),
  @('012521', '英大稳固增强核心一年持有混合A', '0.75', '27.71', '1.96', '0.0147', 4),
  @('001613', '长城久祥灵活配置混合', '0.24', '84.63', '5.98', '0.0144', 1),
  @('014246', '大摩现代服务业混合A', '0.17', '66.96', '8.19', '0.0139', 1),
  @('002292', '诺安益鑫灵活配置混合A', '0.39', '61.16', '3.48', '0.0136', 6),
  @('005569', '中融智选红利股票A', '0.21', '92.04', '6.12', '0.0129', 2),
  @('010077', '湘财长弘灵活配置混合C', '0.31', '68.52', '2.90', '0.0090', 8),
  @('010076', '湘财长弘灵活配置混合A', '0.30', '68.52', '2.90', '0.0087', 8),
  @('001270', '英大灵活配置混合A', '0.28', '93.98', '2.12', '0.0059', 5),
  @('001271', '英大灵活配置混合B', '0.28', '93.98', '2.12', '0.0059', 5),
  @('003446', '英大睿鑫灵活配置混合A', '0.07', '92.71', '7.38', '0.0052', 9),
  @('014247', '大摩现代服务业混合C', '0.06', '66.96', '8.19', '0.0049', 1),
  @('010236', '广发电子信息传媒股票C', '0.13', '89.36', '2.53', '0.0033', 10),
  @('005570', '中融智选红利股票C', '0.03', '92.04', '6.12', '0.0018', 2),
  @('001608', '英大策略优选混合C', '0.02', '91.98', '6.33', '0.0013', 6),
  @('005935', '新疆前海联合润丰灵活配置混合C', '0.01', '87.05', '7.28', '0.0007', 4),
  @('014550', '诺安益鑫灵活配置混合C', '0.02', '61.16', '3.48', '0.0007', 6)
)

$row = 2
foreach ($item in $q3Data) {
    # Column B (fund code) and D:G (scale/position/ratio/marketValue) must
    # stay textual - B has significant leading zeros ("003713") and D:G
    # carry fixed decimal formatting ("7.50") that plain numbers would drop.
    $q3Sheet.Range("B${row}:B${row}").NumberFormat = "@"
    $q3Sheet.Range("D${row}:G${row}").NumberFormat = "@"
    $q3Sheet.Cells.Item($row, 1).Value = $row - 2
    $q3Sheet.Cells.Item($row, 2).Value = $item[0]
    $q3Sheet.Cells.Item($row, 3).Value = $item[1]
    $q3Sheet.Cells.Item($row, 4).Value = $item[2]
    $q3Sheet.Cells.Item($row, 5).Value = $item[3]
    $q3Sheet.Cells.Item($row, 6).Value = $item[4]
    $q3Sheet.Cells.Item($row, 7).Value = $item[5]
    $q3Sheet.Cells.Item($row, 8).Value = $item[6]
    $row++
}

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q3 summary row into "总计", shifting the existing
# rows (2022-Q2 .. 2020-Q4) down by one.
# ---------------------------------------------------------------------------
$lastRow = $totalSheet.UsedRange.Rows.Count

# Shift the existing data rows down by one (read bottom-up so nothing gets
# clobbered before it is copied). Column A (the leading index) is
# renumbered afterwards, so its old value doesn't matter here.
for ($r = $lastRow; $r -ge 2; $r--) {
    $totalSheet.Cells.Item($r + 1, 1).Value = $totalSheet.Cells.Item($r, 1).Value2
    $totalSheet.Cells.Item($r + 1, 2).Value = $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($r + 1, 3).Value = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($r + 1, 4).Value = $totalSheet.Cells.Item($r, 4).Value2
}

# Row 2 (the new row) picks up the same per-column formatting already used
# by every other data row (row 3 still carries it, untouched by the shift
# above): column A keeps the bold/boxed "index" style, B:D stay unstyled.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

# The very last row lost its format when it was first written by the shift
# loop above (it used to be blank), so restore it the same way.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A$($lastRow + 1):A$($lastRow + 1)").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 28
$totalSheet.Cells.Item(2, 4).Value = 0.93

# Renumber the leading index column (A) for every data row so it stays a
# plain 0-based sequence after the insert.
for ($r = 2; $r -le ($lastRow + 1); $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
